# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# For each worker row (16-50) on "Hoja1":
#   - column E "Periodo Mora"   -> replaced with the new (ascending) period list 1705..2003
#   - column F "Salario Basico" -> refreshed minimum-wage figure for that period
#   - column G "Valor Mora"     -> refreshed arrears figure (same for every row)
# and the company logo picture is nudged to the left to match its new anchor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=16; Period="1705"; F=29509; G=781242},
    @{Row=17; Period="1706"; F=29509; G=781242},
    @{Row=18; Period="1707"; F=29509; G=781242},
    @{Row=19; Period="1708"; F=29509; G=781242},
    @{Row=20; Period="1709"; F=29509; G=781242},
    @{Row=21; Period="1710"; F=29509; G=781242},
    @{Row=22; Period="1711"; F=29509; G=781242},
    @{Row=23; Period="1712"; F=29509; G=781242},
    @{Row=24; Period="1801"; F=29509; G=781242},
    @{Row=25; Period="1802"; F=29509; G=781242},
    @{Row=26; Period="1803"; F=29509; G=781242},
    @{Row=27; Period="1804"; F=29509; G=781242},
    @{Row=28; Period="1805"; F=29509; G=781242},
    @{Row=29; Period="1806"; F=29509; G=781242},
    @{Row=30; Period="1807"; F=29509; G=781242},
    @{Row=31; Period="1808"; F=29509; G=781242},
    @{Row=32; Period="1809"; F=31249; G=781242},
    @{Row=33; Period="1810"; F=31249; G=781242},
    @{Row=34; Period="1811"; F=31249; G=781242},
    @{Row=35; Period="1812"; F=31249; G=781242},
    @{Row=36; Period="1901"; F=31249; G=781242},
    @{Row=37; Period="1902"; F=31249; G=781242},
    @{Row=38; Period="1903"; F=31249; G=781242},
    @{Row=39; Period="1904"; F=31249; G=781242},
    @{Row=40; Period="1905"; F=31249; G=781242},
    @{Row=41; Period="1906"; F=31249; G=781242},
    @{Row=42; Period="1907"; F=31249; G=781242},
    @{Row=43; Period="1908"; F=31249; G=781242},
    @{Row=44; Period="1909"; F=31249; G=781242},
    @{Row=45; Period="1910"; F=31249; G=781242},
    @{Row=46; Period="1911"; F=31249; G=781242},
    @{Row=47; Period="1912"; F=31249; G=781242},
    @{Row=48; Period="2001"; F=31249; G=781242},
    @{Row=49; Period="2002"; F=31249; G=781242},
    @{Row=50; Period="2003"; F=31249; G=781242}
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("E$r").Value = $item.Period
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
}

# Reposition the company logo (moved left by 19pt / 241300 EMU; size unchanged)
$shp = $ws.Shapes.Item(1)
$shp.Left = 53.59055118110236
$shp.Width = 76.81889763779527
